# BurnUp chart data update — "SprintProtokol 4 Added"
#
# Sprint 4 (week 7 / row 8) and sprint 5 (week 8 / row 9) results are
# recorded on the "Tabelle1" sheet: the per-sprint work ("Sprint Work(h)",
# column D) is filled in, and the running "Work done (h)" total in column C
# is extended using the same cumulative pattern already used by the rows
# above it (C[n] = C[n-1] + D[n]).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 8 (week 7): 27h worked this sprint.
$ws.Range("D8").Value = 27
$ws.Range("C8").Formula = "=C7+D8"

# Row 9 (week 8): 19h worked this sprint.
$ws.Range("D9").Value = 19
$ws.Range("C9").Formula = "=C8+D9"

# Recalculate so the burn-up chart (which reads Tabelle1!C2:C13) picks up
# the two new cumulative totals (146, 165).
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()
$excel.CalculateFullRebuild()

# Leave the selection where the author ended up after entering the data.
$ws.Range("D16").Select() | Out-Null
